$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title cell (B1) which mirrors the "Nom Code Plaque" value shown in
# the data rows below, so it reflects the new dossier code too.
$ws.Range("B1").Value = "CEMRJ1CO_1013"

# Extend formatting (style s="3") from row 4 down to the new rows 8-11 first,
# so the newly added rows inherit the same borders/style as existing data rows.
$ws.Range("A4:G4").Copy() | Out-Null
$ws.Range("A8:G11").PasteSpecial(-4122) | Out-Null

# New data table for rows 4-11 (header row 3 stays the same: Code IMB, Code plaque,
# Action, ID_ZN, ID_ZE, ID_RGT, Statut code IMB)
$data = @(
    @("IMB/26198/C/02HU", "CEMRJ1CO_1013", "Probation", "SO",              "SO", "SO", ""),
    @("IMB/26198/C/02VU", "CEMRJ1CO_1013", "Probation", "ZN_0319_26_0020", "SO", "SO", "En attente accord Syndic"),
    @("IMB/26198/C/03NR", "CEMRJ1CO_1013", "Probation", "SO",              "SO", "SO", ""),
    @("IMB/26198/C/03YY", "CEMRJ1CO_1013", "Probation", "ZN_0319_26_0007", "SO", "SO", "En attente accord Syndic"),
    @("IMB/26198/C/030S", "CEMRJ1CO_1013", "Probation", "SO",              "SO", "SO", ""),
    @("IMB/26198/C/0374", "CEMRJ1CO_1013", "Probation", "ZN_0319_26_0026", "SO", "SO", "En attente accord Syndic"),
    @("IMB/26198/S/03YW", "CEMRJ1CO_1013", "Probation", "ZN_0319_26_0007", "SO", "SO", "En attente accord Syndic"),
    @("IMB/26198/S/03YX", "CEMRJ1CO_1013", "Probation", "ZN_0319_26_0007", "SO", "SO", "En attente accord Syndic")
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 6).Value = $rowVals[5]
    $ws.Cells.Item($r, 7).Value = $rowVals[6]
}
